$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.257.55'
$ws.Range('E2').Value = '  -0.35%  '
$ws.Range('D3').Value = '3.415.17'
$ws.Range('E3').Value = '  -1.32%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.66'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '177.92'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.46%  '
$ws.Range('E7').Value = '  +4.36%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = '3.416.19'
$ws.Range('E9').Value = '  -1.32%  '
$ws.Range('E10').Value = '  +0.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.95'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.90%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.412'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.17%  '
$ws.Range('D13').Value = '4.011.01'
$ws.Range('E13').Value = '  -1.32%  '
$ws.Range('E14').Value = '  +0.61%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '29.11'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.24%  '
$ws.Range('D16').Value = '66.276.60'
$ws.Range('E16').Value = '  -0.40%  '
$ws.Range('E17').Value = '  +0.43%  '
$ws.Range('D18').Value = '3.402.20'
$ws.Range('E18').Value = '  -1.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.89'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.78'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '367.75'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.55'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '73.20'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.34%  '
$ws.Range('B24').Value = 'PEPE'
$ws.Range('C24').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000127'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.36%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.997'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.31%  '
$ws.Range('B26').Value = 'Polygon'
$ws.Range('C26').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.533'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.76'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.84%  '
$ws.Range('E28').Value = '  +1.25%  '
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.99'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.74'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '23.32'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.98'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.63%  '
$ws.Range('E35').Value = '  -4.84%  '
$ws.Range('E36').Value = '  -1.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '163.76'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.868'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.52%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '27.64'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.79'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.58'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.55%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.43'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.23%  '
$ws.Range('D43').Value = '2.710.26'
$ws.Range('E43').Value = '  -0.79%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.28'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.87%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0685'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '24.94'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.35%  '
$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '337.89'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +9.40%  '
$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '39.93'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.70%  '
$ws.Range('E49').Value = '  -2.25%  '
$ws.Range('E50').Value = '  +2.71%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '31.73'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.91%  '
